$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that change identically for row 2 and row 3
$updates = @{
    "D"  = 0.0158
    "G"  = -0.05072992700729927
    "H"  = -0.05072992700729927
    "I"  = -0.1248175182481752
    "J"  = -0.1248175182481752
    "K"  = -4.49
    "L"  = -0.1638686131386861
    "U"  = 0.985
    "V"  = 0.1865530303030303
    "W"  = 0.3301470588235295
    "X"  = 0.4905957124436082
    "Y"  = -0.1604486536200787
    "Z"  = 0.6045362280469508
    "AA" = -0.07545671167593328
    "AB" = 0.0935476619505296
    "AC" = -0.1690043736264629
    "AD" = 45
    "AF" = 45
    "AG" = 44.015
    "AH" = 0.8949880668257756
    "AI" = 1.122194513715711
    "AJ" = 0.892889745410285
    "AK" = 1.125271634922664
    "AL" = 2.19
    "AM" = 2.174
    "AN" = -115.979381443299
    "AO" = -1.561643835616438
    "AP" = -113.4407216494845
    "AQ" = -1.573137074517019
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("${col}2").Value = $value
    $ws.Range("${col}3").Value = $value
}
